$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.352.22"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "1.591.47"
$ws.Range("E4").Value = "  -0.70%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.507"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "1.814.46"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "1.612.43"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "26.344.77"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("E23").Value = "  -4.17%  "
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("D34").Value = "1.309.08"
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.614"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.65%  "
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -13.84%  "
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("E42").Value = "  +3.81%  "
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.763"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.78%  "
$ws.Range("D46").Value = "1.726.32"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("E47").Value = "  -2.42%  "
$ws.Range("E49").Value = "  -4.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0982"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.42%  "
$ws.Range("E51").Value = "  -1.36%  "
